$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Shift the existing 100-iteration confusion-matrix table down by one
# row (old rows 1-5 -> new rows 2-6) by inserting a blank row at the
# top, then add the "100 Iterations" title above it.
# ---------------------------------------------------------------------
$ws.Rows("1").Insert()
$ws.Range("A1").Value = "100 Iterations"

# ---------------------------------------------------------------------
# Fill in the results (all counts came out 5 of 5) for the 100-
# iteration confusion matrix, now living in rows 4-6.
# ---------------------------------------------------------------------
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 5
$ws.Range("E4").Value = 5
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = 5
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = 5
$ws.Range("E6").Value = 5

# ---------------------------------------------------------------------
# 200-iteration section title + accuracy formula.
# ---------------------------------------------------------------------
$ws.Range("A7").Value = "200 Iterations"
$ws.Range("F7").Formula = "=SUM(C4:E6)/87"

# ---------------------------------------------------------------------
# Build the second confusion-matrix table (rows 8-12). Merge the cells
# first (while they are still blank/default-styled) and only then copy
# the formatting over from the first table, otherwise merging a
# uniformly-bordered range causes Excel to split the border into
# separate top/bottom styles.
# ---------------------------------------------------------------------
$ws.Range("C8:C9").Merge()
$ws.Range("D8:D9").Merge()
$ws.Range("E8:E9").Merge()
$ws.Range("A10:B10").Merge()
$ws.Range("A11:B11").Merge()
$ws.Range("A12:B12").Merge()

$ws.Range("C2:E3").Copy()
$ws.Range("C8:E9").PasteSpecial(-4122)

$ws.Range("A4:E6").Copy()
$ws.Range("A10:E12").PasteSpecial(-4122)

$excel.CutCopyMode = 0

$ws.Range("B8").Value = "Control"
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 2

$ws.Range("A9").Value = "Neither"

$ws.Range("A10").Value = 0
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 5
$ws.Range("E10").Value = 5

$ws.Range("A11").Value = 1
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 5
$ws.Range("E11").Value = 5

$ws.Range("A12").Value = 2
$ws.Range("C12").Value = 4
$ws.Range("D12").Value = 5
$ws.Range("E12").Value = 5

# Accuracy formula for the 200-iteration block.
$ws.Range("F13").Formula = "=SUM(C10:E12)/87"

$ws.Range("D15").Select()
